$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).NumberFormat = "General"
}

# Row 2
$ws.Range("D2").Value = "42.943.32"
$ws.Range("E2").Value = "  -0.34%  "

# Row 3
$ws.Range("D3").Value = "2.303.52"
$ws.Range("E3").Value = "  -0.45%  "

# Row 4
Set-TextValue "D4" "1.00"
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
Set-TextValue "D5" "300.14"
$ws.Range("E5").Value = "  -0.65%  "

# Row 6
Set-TextValue "D6" "97.18"
$ws.Range("E6").Value = "  -2.02%  "

# Row 7
Set-TextValue "D7" "0.512"
$ws.Range("E7").Value = "  +0.50%  "

# Row 8
$ws.Range("E8").Value = "  +0.01%  "

# Row 9
$ws.Range("E9").Value = "  -3.57%  "

# Row 10
Set-TextValue "D10" "35.75"
$ws.Range("E10").Value = "  -0.24%  "

# Row 11
Set-TextValue "D11" "0.0788"
$ws.Range("E11").Value = "  -0.31%  "

# Row 12
Set-TextValue "D12" "17.92"
$ws.Range("E12").Value = "  -0.35%  "

# Row 13
$ws.Range("E13").Value = "  +0.71%  "

# Row 14
$ws.Range("E14").Value = "  -2.36%  "

# Row 15
$ws.Range("D15").Value = "2.661.09"
$ws.Range("E15").Value = "  -0.55%  "

# Row 16
$ws.Range("D16").Value = "2.304.88"
$ws.Range("E16").Value = "  -0.40%  "

# Row 17
Set-TextValue "D17" "0.778"
$ws.Range("E17").Value = "  -1.61%  "

# Row 18
$ws.Range("D18").Value = "42.894.70"
$ws.Range("E18").Value = "  -0.29%  "

# Row 19
Set-TextValue "D19" "12.73"
$ws.Range("E19").Value = "  -5.50%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0904"
$ws.Range("E20").Value = "  -0.83%  "

# Row 21
Set-TextValue "D21" "6.04"
$ws.Range("E21").Value = "  -2.21%  "

# Row 22
Set-TextValue "D22" "67.87"
$ws.Range("E22").Value = "  -0.27%  "

# Row 23
Set-TextValue "D23" "240.19"
$ws.Range("E23").Value = "  -0.19%  "

# Row 24
Set-TextValue "D24" "2.14"
$ws.Range("E24").Value = "  -1.26%  "

# Row 25
$ws.Range("E25").Value = "  +0.05%  "

# Row 26
Set-TextValue "D26" "2.42"
$ws.Range("E26").Value = "  -0.85%  "

# Row 27
$ws.Range("E27").Value = "  -0.08%  "

# Row 28
Set-TextValue "D28" "25.41"
$ws.Range("E28").Value = "  +1.71%  "

# Row 29
Set-TextValue "D29" "165.74"
$ws.Range("E29").Value = "  -1.50%  "

# Row 30
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D30" "2.03"
$ws.Range("E30").Value = "  -1.19%  "

# Row 31
$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D31" "9.04"
$ws.Range("E31").Value = "  -1.70%  "

# Row 32
Set-TextValue "D32" "33.04"
$ws.Range("E32").Value = "  -1.48%  "

# Row 33
Set-TextValue "D33" "4.90"
$ws.Range("E33").Value = "  -0.53%  "

# Row 34
$ws.Range("E34").Value = "  +0.09%  "

# Row 35
Set-TextValue "D35" "5.02"
$ws.Range("E35").Value = "  -4.14%  "

# Row 36
Set-TextValue "D36" "16.96"
$ws.Range("E36").Value = "  -7.53%  "

# Row 37
$ws.Range("E37").Value = "  -1.19%  "

# Row 38
Set-TextValue "D38" "0.0686"
$ws.Range("E38").Value = "  -1.15%  "

# Row 39
$ws.Range("E39").Value = "  -1.34%  "

# Row 40
$ws.Range("E40").Value = "  -2.71%  "

# Row 41
$ws.Range("E41").Value = "  -0.15%  "

# Row 42
$ws.Range("E42").Value = "  -1.42%  "

# Row 43
$ws.Range("D43").Value = "2.013.28"
$ws.Range("E43").Value = "  +0.75%  "

# Row 44
Set-TextValue "D44" "0.0281"
$ws.Range("E44").Value = "  -2.71%  "

# Row 45
Set-TextValue "D45" "10.17"
$ws.Range("E45").Value = "  +0.59%  "

# Row 46
Set-TextValue "D46" "2.12"
$ws.Range("E46").Value = "  -2.31%  "

# Row 47
Set-TextValue "D47" "17.37"

# Row 48
Set-TextValue "D48" "2.78"
$ws.Range("E48").Value = "  -2.03%  "

# Row 49
Set-TextValue "D49" "2.92"
$ws.Range("E49").Value = "  -2.35%  "

# Row 50
Set-TextValue "D50" "53.56"
$ws.Range("E50").Value = "  -2.53%  "

# Row 51
$ws.Range("D51").Value = "2.526.73"
$ws.Range("E51").Value = "  -0.58%  "

